# [Fonds de solidarite] Add 2020-08-28 data
# Updates nombre_aides (col C) and montant_total (col D) for the rows whose
# underlying data changed between the 2020-08-27 and 2020-08-28 exports.
#
# Values in these columns are stored as text (e.g. "1021", not the number
# 1021), so the NumberFormat/Style dance below forces Excel to keep the new
# value as text instead of silently re-typing it as a number, while leaving
# the cell's final style untouched (back to the default "Normal" style).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

# row -> (nombre_aides, montant_total)
Set-TextValue "C3"  "1026"
Set-TextValue "D3"  "3286820.33"

Set-TextValue "C4"  "423"
Set-TextValue "D4"  "1754198.25"

Set-TextValue "C6"  "30"
Set-TextValue "D6"  "196643.82"

Set-TextValue "C10" "365"
Set-TextValue "D10" "1311761.71"

Set-TextValue "C34" "577"
Set-TextValue "D34" "1903771.66"

Set-TextValue "C35" "234"
Set-TextValue "D35" "1192368.11"

Set-TextValue "C36" "76"
Set-TextValue "D36" "427894.00"

Set-TextValue "C38" "23"
Set-TextValue "D38" "50200.00"

Set-TextValue "C52" "601"
Set-TextValue "D52" "2135615.21"

Set-TextValue "C53" "266"
Set-TextValue "D53" "1182878.76"

Set-TextValue "C55" "26"
Set-TextValue "D55" "143213.00"

Set-TextValue "C56" "25"
Set-TextValue "D56" "80220.65"

Set-TextValue "C83" "897"
Set-TextValue "D83" "2880909.26"
